$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.839.59'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '2.349.77'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '544.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.07%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.525'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.13%  '
$ws.Range('D9').Value = '2.347.91'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.343'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '24.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.16%  '
$ws.Range('D15').Value = '2.772.52'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '61.025.63'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '2.350.47'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.37%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '319.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.66%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.78%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.95'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '496.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('D31').Value = '0.0₃0860'
$ws.Range('E31').Value = '  -6.07%  '
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('E34').Value = '  -2.00%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.27%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.24'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.99%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '144.71'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.19%  '
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '143.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.06%  '
$ws.Range('E45').Value = '  +1.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0515'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('E48').Value = '  -4.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.568'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0220'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.66%  '
